$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 22 (this shifts existing rows 22-44 down to 23-45)
$ws.Rows.Item(22).Insert()

# Fill the new row 22 with the new test step
$ws.Range("A22").Value = 'Der Nutzer wählt für Spieler 1 "Mensch", für Spieler 2 "KI 1" und klickt auf das Feld "Spiel starten".'
$ws.Range("B22").Value = "Es wird in die Spielansicht gewechselt. Das Spielfeld ist leer. Der Graph zeigt das leere Feld und alle Äquivalenzklassenvertreter der möglichen Folgezustände an."

# The insert carried over the bold header style from row 21; reset the new
# row back to the default (non-bold) style used by ordinary data rows.
$ws.Range("A22:B22").Font.Bold = $false

# Move selection to match the final state
$ws.Range("B21").Select()
